# Applies the "updated referee response to indicate done things" edit:
#  1. Merge the split run "We have ad" + "ded this figure (now Figure X)."
#     into a single run "We have added this figure (now Figure X)."
#     (keeping the _GoBack bookmark in place).
#  2. Recolor several "done" responses from grey/red to green (92D050):
#     - "We have included a table of chemical abundances measured from the
#        Magellan/MIKE spectra." (run only)
#     - "Good catch! We have included these in the updated manuscript."
#        (paragraph mark + run)
#     - "We have expanded the discussion on this point. ... etc ... we
#        comment on it for completeness." (paragraph mark + its 3 runs)

$d = $word.ActiveDocument

# Word BGR-packed color value for RGB hex 92D050 (R=0x92,G=0xD0,B=0x50)
$green = 0x50 * 65536 + 0xD0 * 256 + 0x92

# ---------------------------------------------------------------------
# 1) Merge the "We have ad" / "ded this figure (now Figure X)." runs.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    # Range.Text includes the trailing paragraph-mark character, so trim it
    # before comparing against the literal target text.
    $paraText = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($paraText -eq "We have added this figure (now Figure X).") {
        # There are two paragraphs in the document with this same final
        # text; only the one still split across two runs (identifiable by
        # the presence of the lone "_GoBack" bookmark) needs editing.
        $isSplit = $p.Range.WordOpenXML -like "*_GoBack*"
        if ($isSplit) {
            # Delete the second (redundant) run's text first, while the
            # phrase is still unambiguous.
            $delRange = $p.Range
            $delRange.Find.Execute("ded this figure (now Figure X).", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

            # Now merge the remaining text into the first run.
            $mergeRange = $p.Range
            $mergeRange.Find.Execute("We have ad", $true, $false, $false, $false, $false, $true, 1, $false, "We have added this figure (now Figure X).", 2) | Out-Null
        }
    }
}

# ---------------------------------------------------------------------
# 2) Recolor the three "done" responses to green (92D050).
# ---------------------------------------------------------------------

# a) Only the first run of this paragraph changes color (the rest of the
#    paragraph, "- A table with ...", stays as-is).
$r = $d.Content
$r.Find.Execute("We have included a table of chemical abundances measured from the Magellan/MIKE spectra.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Font.Color = $green

# b) Whole paragraph (paragraph mark + its single run) changes color.
foreach ($i in 1..$paras.Count) {
    $p = $paras.Item($i)
    $paraText = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($paraText -eq "Good catch! We have included these in the updated manuscript.") {
        $p.Range.Font.Color = $green
    }
}

# c) Whole paragraph (paragraph mark + its three runs) changes color; the
#    following paragraph ("This claim is elaborated on...") is untouched.
foreach ($i in 1..$paras.Count) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "We have expanded the discussion on this point.*completeness.*") {
        $p.Range.Font.Color = $green
    }
}
